# Updates cryptos list values (prices / 1h volume %) per the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($cellRef, $val) {
    $rng = $ws.Range($cellRef)
    # Force text storage so numeric-looking strings (e.g. "0.9579")
    # are not silently reinterpreted as numbers by Excel.
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = "Normal"
}

$ws.Range("D2").Value = '20.406.39'
$ws.Range("E2").Value = '  +2.15%  '
$ws.Range("D3").Value = '1.464.08'
$ws.Range("E3").Value = '  +3.11%  '
$ws.Range("E4").Value = '  +0.72%  '
Set-TextCell "D5" '0.9579'
$ws.Range("E5").Value = '  -4.20%  '
Set-TextCell "D6" '274.89'
$ws.Range("E6").Value = '  -0.65%  '
Set-TextCell "D7" '0.3647'
$ws.Range("E7").Value = '  -1.06%  '
Set-TextCell "D8" '0.3061'
$ws.Range("E8").Value = '  -1.31%  '
Set-TextCell "D9" '39.82'
$ws.Range("E9").Value = '  -0.03%  '
Set-TextCell "D10" '1.041'
$ws.Range("E10").Value = '  -0.49%  '
Set-TextCell "D11" '0.06596'
$ws.Range("E11").Value = '  +0.71%  '
$ws.Range("E12").Value = '  -0.02%  '
Set-TextCell "D13" '18.16'
$ws.Range("E13").Value = '  +2.55%  '
Set-TextCell "D14" '5.410'
$ws.Range("E14").Value = '  -2.04%  '
Set-TextCell "D15" '6.142'
$ws.Range("E15").Value = '  -1.27%  '
Set-TextCell "D16" '0.00001023'
$ws.Range("E16").Value = '  -0.05%  '
$ws.Range("D17").Value = '1.464.98'
$ws.Range("E17").Value = '  +2.97%  '
Set-TextCell "D18" '0.9735'
$ws.Range("E18").Value = '  -2.64%  '
Set-TextCell "D19" '0.05875'
$ws.Range("E19").Value = '  +3.30%  '
Set-TextCell "D20" '69.34'
$ws.Range("E20").Value = '  -3.00%  '
Set-TextCell "D21" '5.431'
$ws.Range("E21").Value = '  -3.53%  '
Set-TextCell "D22" '14.31'
$ws.Range("E22").Value = '  -3.15%  '
Set-TextCell "D23" '10.91'
$ws.Range("E23").Value = '  -0.93%  '
Set-TextCell "D24" '2.242'
$ws.Range("E24").Value = '  +0.11%  '
$ws.Range("D25").Value = '20.427.90'
$ws.Range("E25").Value = '  +2.07%  '
Set-TextCell "D26" '142.86'
$ws.Range("E26").Value = '  +7.21%  '
Set-TextCell "D27" '2.075'
$ws.Range("E27").Value = '  -9.54%  '
$ws.Range("E28").Value = '  -1.27%  '
$ws.Range("D29").Value = '1.620.03'
$ws.Range("E29").Value = '  +2.42%  '
Set-TextCell "D30" '114.05'
$ws.Range("E30").Value = '  +3.36%  '
Set-TextCell "D31" '3.842'
$ws.Range("E31").Value = '  -2.18%  '
Set-TextCell "D32" '4.946'
$ws.Range("E32").Value = '  -6.25%  '
Set-TextCell "D33" '0.07884'
$ws.Range("E33").Value = '  +1.62%  '
Set-TextCell "D34" '0.7919'
$ws.Range("E34").Value = '  -3.92%  '
Set-TextCell "D35" '1.503'
$ws.Range("E35").Value = '  +0.87%  '
Set-TextCell "D36" '0.05719'
$ws.Range("E36").Value = '  -2.28%  '
Set-TextCell "D37" '1.146'
$ws.Range("E37").Value = '  +3.85%  '
Set-TextCell "D38" '4.690'
$ws.Range("E38").Value = '  -4.90%  '
$ws.Range("B39").Value = 'FraxShare'
$ws.Range("C39").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
Set-TextCell "D39" '7.682'
$ws.Range("E39").Value = '  -7.45%  '
$ws.Range("B40").Value = 'VeChain'
$ws.Range("C40").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
Set-TextCell "D40" '0.02034'
$ws.Range("E40").Value = '  -1.71%  '
$ws.Range("B41").Value = 'Frax'
$ws.Range("C41").Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
Set-TextCell "D41" '0.9605'
$ws.Range("E41").Value = '  -3.81%  '
$ws.Range("B42").Value = 'Aptos'
$ws.Range("C42").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
Set-TextCell "D42" '10.34'
$ws.Range("E42").Value = '  -1.75%  '
Set-TextCell "D43" '0.1856'
$ws.Range("E43").Value = '  -1.80%  '
Set-TextCell "D44" '0.5271'
$ws.Range("E44").Value = '  -1.16%  '
Set-TextCell "D45" '3.483'
$ws.Range("E45").Value = '  -1.69%  '
Set-TextCell "D46" '11.97'
$ws.Range("E46").Value = '  -4.26%  '
Set-TextCell "D47" '116.86'
$ws.Range("E47").Value = '  +0.10%  '
Set-TextCell "D48" '0.5171'
$ws.Range("E48").Value = '  -0.78%  '
Set-TextCell "D49" '1.753'
$ws.Range("E49").Value = '  -1.44%  '
Set-TextCell "D50" '0.06421'
$ws.Range("E50").Value = '  +3.62%  '
Set-TextCell "D51" '0.9918'
$ws.Range("E51").Value = '  -0.67%  '
